# Update the "last_edited_time" (column D) values on the
# LUY_KE_NGAY_CAN_THO sheet to reflect the new Notion sync timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D9").Value   = "2024-07-04T09:36:00.000Z"
$ws.Range("D10:D30").Value = "2024-07-04T09:33:00.000Z"
$ws.Range("D31:D76").Value = "2024-07-04T09:34:00.000Z"
$ws.Range("D77:D115").Value = "2024-07-04T09:35:00.000Z"
